$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update row 14 (E0XY instruction): overflow -> not overflow
$ws.Range("C14").Value = "jnc mem"
$ws.Range("B14").Value = "Jump if not overflow"
$ws.Range("D14").Value = "If (!overflow) goto mem"

# Update the active selection to D14
$ws.Range("D14").Select()
